# Final commit with 4 testcases
# Update the second test-data row (username / email / password) and
# move the active selection from D2 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update password first so the shared-string table ends up ordered the
# same way as in the target workbook: Aldod@334; keeps its slot, then
# the new username/email strings are appended after it.
$ws.Range("C2").Value = "Aldod@334;"
$ws.Range("A2").Value = "TesdgsLE"
$ws.Range("B2").Value = "testersvbsksle@xyz.com"

# Move the selected/active cell from D2 to B2.
[void]$ws.Range("B2").Select()
